$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6923383333333333
$ws.Range("H2").Value = 2.077015
$ws.Range("I2").Value = 0.2954746751438028
$ws.Range("J2").Value = 0.2954746751438028
$ws.Range("M2").Value = 0.2949276666666666
$ws.Range("N2").Value = 0.884783
$ws.Range("O2").Value = 0.01958234361069057
$ws.Range("P2").Value = 0.01958234361069057
$ws.Range("Q2").Value = 0.2041897291938889
$ws.Range("R2").Value = 1.837707562745
$ws.Range("S2").Value = 0.005786086616923118
$ws.Range("T2").Value = 0.005786086616923118
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6923383333333333
$ws.Range("H3").Value = 2.077015
$ws.Range("I3").Value = 0.2954746751438028
$ws.Range("J3").Value = 0.2954746751438028
$ws.Range("O3").Value = 0.4183293461133303
$ws.Range("P3").Value = 0.4183293461133303
$ws.Range("Q3").Value = 4.362019051187778
$ws.Range("R3").Value = 39.25817146069
$ws.Range("S3").Value = 0.1236057276459557
$ws.Range("T3").Value = 0.1236057276459557
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6923383333333333
$ws.Range("H4").Value = 2.077015
$ws.Range("I4").Value = 0.2954746751438028
$ws.Range("J4").Value = 0.2954746751438028
$ws.Range("M4").Value = 0.7642679999999999
$ws.Range("N4").Value = 2.292804
$ws.Range("O4").Value = 0.05074518357604722
$ws.Range("P4").Value = 0.05074518357604721
$ws.Range("Q4").Value = 0.5291320333399999
$ws.Range("R4").Value = 4.762188300059999
$ws.Range("S4").Value = 0.01499391663224519
$ws.Range("T4").Value = 0.01499391663224518
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6923383333333333
$ws.Range("H5").Value = 2.077015
$ws.Range("I5").Value = 0.2954746751438028
$ws.Range("J5").Value = 0.2954746751438028
$ws.Range("M5").Value = 7.701286333333333
$ws.Range("N5").Value = 23.103859
$ws.Range("O5").Value = 0.5113431266999319
$ws.Range("P5").Value = 0.5113431266999319
$ws.Range("Q5").Value = 5.331895744542777
$ws.Range("R5").Value = 47.98706170088499
$ws.Range("S5").Value = 0.1510889442486787
$ws.Range("T5").Value = 0.1510889442486787
$ws.Range("I6").Value = 0.4065337983884298
$ws.Range("J6").Value = 0.4065337983884298
$ws.Range("M6").Value = 0.2949276666666666
$ws.Range("N6").Value = 0.884783
$ws.Range("O6").Value = 0.01958234361069057
$ws.Range("P6").Value = 0.01958234361069057
$ws.Range("Q6").Value = 0.2809378711075555
$ws.Range("R6").Value = 2.528440839968
$ws.Range("S6").Value = 0.007960884529401437
$ws.Range("T6").Value = 0.007960884529401437
$ws.Range("I7").Value = 0.4065337983884298
$ws.Range("J7").Value = 0.4065337983884298
$ws.Range("O7").Value = 0.4183293461133303
$ws.Range("P7").Value = 0.4183293461133303
$ws.Range("S7").Value = 0.1700650180528003
$ws.Range("T7").Value = 0.1700650180528003
$ws.Range("I8").Value = 0.4065337983884298
$ws.Range("J8").Value = 0.4065337983884298
$ws.Range("M8").Value = 0.7642679999999999
$ws.Range("N8").Value = 2.292804
$ws.Range("O8").Value = 0.05074518357604722
$ws.Range("P8").Value = 0.05074518357604721
$ws.Range("Q8").Value = 0.7280152021759999
$ws.Range("R8").Value = 6.552136819583999
$ws.Range("S8").Value = 0.02062963222908864
$ws.Range("T8").Value = 0.02062963222908864
$ws.Range("I9").Value = 0.4065337983884298
$ws.Range("J9").Value = 0.4065337983884298
$ws.Range("M9").Value = 7.701286333333333
$ws.Range("N9").Value = 23.103859
$ws.Range("O9").Value = 0.5113431266999319
$ws.Range("P9").Value = 0.5113431266999319
$ws.Range("Q9").Value = 7.33597838320711
$ws.Range("R9").Value = 66.023805448864
$ws.Range("S9").Value = 0.2078782635771395
$ws.Range("T9").Value = 0.2078782635771395
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.3086113333333333
$ws.Range("H10").Value = 0.925834
$ws.Range("I10").Value = 0.1317084856811759
$ws.Range("J10").Value = 0.1317084856811759
$ws.Range("M10").Value = 0.2949276666666666
$ws.Range("N10").Value = 0.884783
$ws.Range("O10").Value = 0.01958234361069057
$ws.Range("P10").Value = 0.01958234361069057
$ws.Range("Q10").Value = 0.09101802044688889
$ws.Range("R10").Value = 0.819162184022
$ws.Range("S10").Value = 0.002579160823052505
$ws.Range("T10").Value = 0.002579160823052505
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.3086113333333333
$ws.Range("H11").Value = 0.925834
$ws.Range("I11").Value = 0.1317084856811759
$ws.Range("J11").Value = 0.1317084856811759
$ws.Range("O11").Value = 0.4183293461133303
$ws.Range("P11").Value = 0.4183293461133303
$ws.Range("Q11").Value = 1.944379576573778
$ws.Range("R11").Value = 17.499416189164
$ws.Range("S11").Value = 0.05509752469258324
$ws.Range("T11").Value = 0.05509752469258323
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.3086113333333333
$ws.Range("H12").Value = 0.925834
$ws.Range("I12").Value = 0.1317084856811759
$ws.Range("J12").Value = 0.1317084856811759
$ws.Range("M12").Value = 0.7642679999999999
$ws.Range("N12").Value = 2.292804
$ws.Range("O12").Value = 0.05074518357604722
$ws.Range("P12").Value = 0.05074518357604721
$ws.Range("Q12").Value = 0.235861766504
$ws.Range("R12").Value = 2.122755898536
$ws.Range("S12").Value = 0.006683571284414457
$ws.Range("T12").Value = 0.006683571284414456
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.3086113333333333
$ws.Range("H13").Value = 0.925834
$ws.Range("I13").Value = 0.1317084856811759
$ws.Range("J13").Value = 0.1317084856811759
$ws.Range("M13").Value = 7.701286333333333
$ws.Range("N13").Value = 23.103859
$ws.Range("O13").Value = 0.5113431266999319
$ws.Range("P13").Value = 0.5113431266999319
$ws.Range("Q13").Value = 2.376704243711778
$ws.Range("R13").Value = 21.390338193406
$ws.Range("S13").Value = 0.06734822888112568
$ws.Range("T13").Value = 0.06734822888112568
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.3896243333333334
$ws.Range("H14").Value = 1.168873
$ws.Range("I14").Value = 0.1662830407865915
$ws.Range("J14").Value = 0.1662830407865915
$ws.Range("M14").Value = 0.2949276666666666
$ws.Range("N14").Value = 0.884783
$ws.Range("O14").Value = 0.01958234361069057
$ws.Range("P14").Value = 0.01958234361069057
$ws.Range("Q14").Value = 0.1149109955065555
$ws.Range("R14").Value = 1.034198959559
$ws.Range("S14").Value = 0.003256211641313509
$ws.Range("T14").Value = 0.003256211641313509
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.3896243333333334
$ws.Range("H15").Value = 1.168873
$ws.Range("I15").Value = 0.1662830407865915
$ws.Range("J15").Value = 0.1662830407865915
$ws.Range("O15").Value = 0.4183293461133303
$ws.Range("P15").Value = 0.4183293461133303
$ws.Range("Q15").Value = 2.454795123973111
$ws.Range("R15").Value = 22.093156115758
$ws.Range("S15").Value = 0.06956107572199104
$ws.Range("T15").Value = 0.06956107572199104
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.3896243333333334
$ws.Range("H16").Value = 1.168873
$ws.Range("I16").Value = 0.1662830407865915
$ws.Range("J16").Value = 0.1662830407865915
$ws.Range("M16").Value = 0.7642679999999999
$ws.Range("N16").Value = 2.292804
$ws.Range("O16").Value = 0.05074518357604722
$ws.Range("P16").Value = 0.05074518357604721
$ws.Range("Q16").Value = 0.297777409988
$ws.Range("R16").Value = 2.679996689892
$ws.Range("S16").Value = 0.008438063430298929
$ws.Range("T16").Value = 0.008438063430298929
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.3896243333333334
$ws.Range("H17").Value = 1.168873
$ws.Range("I17").Value = 0.1662830407865915
$ws.Range("J17").Value = 0.1662830407865915
$ws.Range("M17").Value = 7.701286333333333
$ws.Range("N17").Value = 23.103859
$ws.Range("O17").Value = 0.5113431266999319
$ws.Range("P17").Value = 0.5113431266999319
$ws.Range("Q17").Value = 3.000608553434111
$ws.Range("R17").Value = 27.005476980907
$ws.Range("S17").Value = 0.08502768999298797
$ws.Range("T17").Value = 0.08502768999298797
